# Swap columns C and D (codeforiati:group-name <-> codeforiati:group-code)
# for the header row and every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range's last row reliably
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
